# Applies the "rough draft of paper" edit:
#  - Renames Sheet1 -> "Survey Responses"
#  - Adds a new worksheet "Topic Subjects" after it (becomes the active tab)
#  - Re-labels Survey Responses header row: A1 = RespondentID, B1:K1 = Topic1..Topic10
#  - Narrows column A on Survey Responses and selects B1:K1
#  - Populates Topic Subjects with a topic/subject lookup table

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Survey Responses"

# Header row: write the Topic columns first, then the RespondentID cell,
# so the shared-string table is built up in the same order as the target.
$ws1.Range("B1").Value = "Topic1"
$ws1.Range("C1").Value = "Topic2"
$ws1.Range("D1").Value = "Topic3"
$ws1.Range("E1").Value = "Topic4"
$ws1.Range("F1").Value = "Topic5"
$ws1.Range("G1").Value = "Topic6"
$ws1.Range("H1").Value = "Topic7"
$ws1.Range("I1").Value = "Topic8"
$ws1.Range("J1").Value = "Topic9"
$ws1.Range("K1").Value = "Topic10"
$ws1.Range("A1").Value = "RespondentID"

# Narrow the RespondentID column and mark B1:K1 as the active selection.
$ws1.Columns.Item(1).ColumnWidth = 10.166666666666666
$ws1.Range("B1:K1").Select() | Out-Null

# Add the new "Topic Subjects" sheet right after "Survey Responses".
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Topic Subjects"

$subjects = @(
  "Availability during production failures",
  "Financial incentives in the supply chain",
  "Assessing competing products",
  "Safety and effectiveness",
  "Co-ordination of government financing programs",
  "Payment for vaccination and follow up care",
  "Side Effect monitoring",
  "Bioterrorism ",
  "Communicating Risks and improving knowledge",
  "Global Health"
)

for ($i = 0; $i -lt 10; $i++) {
  $row = $i + 2
  $ws2.Cells.Item($row, 1).Value = $ws1.Cells.Item(1, $i + 2).Value2
  $ws2.Cells.Item($row, 2).Value = $subjects[$i]
}

$ws2.Range("A1").Value = "topic"
$ws2.Range("B1").Value = "subject"

$ws2.Columns.Item(2).ColumnWidth = 35.166666666666664
$ws2.Range("C1:C11").Select() | Out-Null
